# Edit full watchlist: remove commodity/crypto symbols (USOIL, SILVER, GOLD,
# BTCUSD, LINKUSD) from the top of the list and append the new tickers for
# the watchlist set starting 20.03.2022.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows 3-7 (USOIL, SILVER, GOLD, BTCUSD, LINKUSD), shifting cells up.
$ws.Range("A3:A7").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp)

# New tickers to append to the bottom of the watchlist.
$newSymbols = @(
  "CAN","HES","MOS","DWAC","OKE","NTR","CMC","WHD","ORLY","GFI",
  "RS","OXY","ANTM","REGN","COST","JBHT","UNH","HRMY","IRM","ODFL",
  "BLDR","AVGO","CVS","GFS","SHELL","BRO","MAR","JNPR","MET","XOM",
  "DKNG","COP","V","FDX","KBH","TMST","SSRM","ACLS"
)

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End([Microsoft.Office.Interop.Excel.XlDirection]::xlUp).Row

for ($i = 0; $i -lt $newSymbols.Length; $i++) {
    $ws.Cells.Item($lastRow + 1 + $i, 1).Value2 = $newSymbols[$i]
}

# Update the selection / view to reflect the new extent of the list, mirroring
# the workbook being scrolled to the end of the (now longer) watchlist.
$lastDataRow = $lastRow + $newSymbols.Length
$ws.Application.ActiveWindow.ScrollRow = 123
$ws.Range("A" + ($lastDataRow + 1)).Select()

$wb.Save()
